$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3955.9546
$ws.Range("I76").Value = 3160
$ws.Range("J76").Value = 4619.25
$ws.Range("K76").Value = 3160
$ws.Range("L76").Value = 4619.25
$ws.Range("M76").Value = -2845
$ws.Range("N76").Value = -5249.25
$ws.Range("H79").Value = 3955.9546
$ws.Range("I79").Value = 3160
$ws.Range("J79").Value = 4619.25
$ws.Range("K79").Value = 3160
$ws.Range("L79").Value = 4619.25
$ws.Range("M79").Value = -2068
$ws.Range("N79").Value = -6803.25
$ws.Range("H113").Value = 2012.4
$ws.Range("I113").Value = 1829.6
$ws.Range("J113").Value = 2378
$ws.Range("K113").Value = 1829.6
$ws.Range("L113").Value = 2378
$ws.Range("M113").Value = 1424.4
$ws.Range("N113").Value = -8886
$ws.Range("H116").Value = 2137.0908
$ws.Range("I116").Value = 2163.5
$ws.Range("K116").Value = 2163.5
$ws.Range("M116").Value = 1278.5
$ws.Range("H132").Value = 3386.6904
$ws.Range("I132").Value = 2664.0286
$ws.Range("K132").Value = 7992.085800000001
$ws.Range("M132").Value = -5462.085800000001
$ws.Range("H133").Value = 47000
$ws.Range("J133").Value = 47000
$ws.Range("L133").Value = 47000
$ws.Range("N133").Value = -57120
$ws.Range("H135").Value = 15152919
$ws.Range("I135").Value = 1472.4333
$ws.Range("J135").Value = 166667400
$ws.Range("K135").Value = 13251.8997
$ws.Range("L135").Value = 1500006600
$ws.Range("M135").Value = -10716.8997
$ws.Range("N135").Value = -1500011670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14287751
$ws.Range("I61").Value = 15153514
$ws.Range("J61").Value = 2649
$ws.Range("K61").Value = 15153514
$ws.Range("L61").Value = 2649
$ws.Range("M61").Value = -15153302
$ws.Range("N61").Value = -3073
$ws.Range("H102").Value = 2095.5715
$ws.Range("I102").Value = 1804
$ws.Range("J102").Value = 2824.5
$ws.Range("K102").Value = 1804
$ws.Range("L102").Value = 2824.5
$ws.Range("M102").Value = -182
$ws.Range("N102").Value = -6068.5
$ws.Range("H132").Value = 10419174
$ws.Range("I132").Value = 13891124
$ws.Range("J132").Value = 3323.3333
$ws.Range("K132").Value = 41673372
$ws.Range("L132").Value = 9969.999899999999
$ws.Range("M132").Value = -41670842
$ws.Range("N132").Value = -15029.9999
$ws.Range("H136").Value = 14287751
$ws.Range("I136").Value = 15153514
$ws.Range("J136").Value = 2649
$ws.Range("K136").Value = 45460542
$ws.Range("L136").Value = 7947
$ws.Range("M136").Value = -45457992
$ws.Range("N136").Value = -13047

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 820.6111
$ws.Range("I99").Value = 775.38464
$ws.Range("J99").Value = 938.2
$ws.Range("K99").Value = 775.38464
$ws.Range("L99").Value = 938.2
$ws.Range("M99").Value = 722.61536
$ws.Range("N99").Value = -3934.2
$ws.Range("H105").Value = 4435.84
$ws.Range("I105").Value = 2849.5
$ws.Range("J105").Value = 4738
$ws.Range("K105").Value = 2849.5
$ws.Range("L105").Value = 4738
$ws.Range("M105").Value = -1102.5
$ws.Range("N105").Value = -8232
$ws.Range("H134").Value = 5711.1904
$ws.Range("I134").Value = 4385.5
$ws.Range("J134").Value = 8362.571
$ws.Range("K134").Value = 13156.5
$ws.Range("L134").Value = 25087.713
$ws.Range("M134").Value = -10621.5
$ws.Range("N134").Value = -30157.713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 16666850
$ws.Range("I12").Value = 16666850
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 16666850
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -16666680
$ws.Range("N12").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3945.4517
$ws.Range("I134").Value = 1995.45
$ws.Range("J134").Value = 7490.909
$ws.Range("K134").Value = 5986.35
$ws.Range("L134").Value = 22472.727
$ws.Range("M134").Value = -916.3500000000004
$ws.Range("N134").Value = -32612.727
$ws.Range("H139").Value = 2768.9473
$ws.Range("I139").Value = 1506.875
$ws.Range("J139").Value = 9500
$ws.Range("K139").Value = 4520.625
$ws.Range("L139").Value = 28500
$ws.Range("M139").Value = 619.375
$ws.Range("N139").Value = -38780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4017.5386
$ws.Range("I102").Value = 4560.3335
$ws.Range("K102").Value = 4560.3335
$ws.Range("M102").Value = -2938.3335
$ws.Range("H113").Value = 78236.46000000001
$ws.Range("I113").Value = 250727.75
$ws.Range("J113").Value = 1573.6666
$ws.Range("K113").Value = 250727.75
$ws.Range("L113").Value = 1573.6666
$ws.Range("M113").Value = -248557.75
$ws.Range("N113").Value = -5913.6666
$ws.Range("H126").Value = 5878.5
$ws.Range("J126").Value = 6141
$ws.Range("L126").Value = 18423
$ws.Range("N126").Value = -23363
$ws.Range("H132").Value = 4694.306
$ws.Range("I132").Value = 4995.9736
$ws.Range("J132").Value = 3652.182
$ws.Range("K132").Value = 14987.9208
$ws.Range("L132").Value = 10956.546
$ws.Range("M132").Value = -12457.9208
$ws.Range("N132").Value = -16016.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2753.75
$ws.Range("I16").Value = 2731.3635
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2731.3635
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -2561.3635
$ws.Range("N16").Value = -3340
$ws.Range("H17").Value = 10000000
$ws.Range("J17").Value = 10000000
$ws.Range("L17").Value = 10000000
$ws.Range("N17").Value = -10000340
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null
$ws.Range("H22").Value = 1380.32
$ws.Range("I22").Value = 612.5
$ws.Range("J22").Value = 1526.5714
$ws.Range("K22").Value = 612.5
$ws.Range("L22").Value = 1526.5714
$ws.Range("M22").Value = -317.5
$ws.Range("N22").Value = -2116.5714
$ws.Range("H27").Value = 1380.32
$ws.Range("I27").Value = 612.5
$ws.Range("J27").Value = 1526.5714
$ws.Range("K27").Value = 612.5
$ws.Range("L27").Value = 1526.5714
$ws.Range("M27").Value = -505.5
$ws.Range("N27").Value = -1740.5714
$ws.Range("H46").Value = 924
$ws.Range("I46").Value = 699.5714
$ws.Range("K46").Value = 699.5714
$ws.Range("M46").Value = -511.5714
$ws.Range("H93").Value = 1366.6666
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 1550
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 1550
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -4046
$ws.Range("H122").Value = 6660.2964
$ws.Range("I122").Value = 7264.25
$ws.Range("J122").Value = 5781.8184
$ws.Range("K122").Value = 21792.75
$ws.Range("L122").Value = 17345.4552
$ws.Range("M122").Value = -19342.75
$ws.Range("N122").Value = -22245.4552
$ws.Range("H136").Value = 62503124
$ws.Range("I136").Value = 62503124
$ws.Range("K136").Value = 187509372
$ws.Range("M136").Value = -187506822

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2650
$ws.Range("I17").Value = 2650
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2650
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2478
$ws.Range("N17").Value = $null
$ws.Range("H41").Value = 6855
$ws.Range("I41").Value = 4999
$ws.Range("J41").Value = 7226.2
$ws.Range("K41").Value = 4999
$ws.Range("L41").Value = 7226.2
$ws.Range("M41").Value = -4609
$ws.Range("N41").Value = -8006.2
